$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" property (row 7) previously had no value; set it to the
# literal text "false" (force text so it isn't auto-coerced to a boolean).
$ws.Cells.Item(7, 2).Value = "'false"

# "Date" property (row 8) value is refreshed to the new publication date.
$ws.Cells.Item(8, 2).Value = "2025-11-04T10:04:56+00:00"
